$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column A ("Match ID"); everything that used to live in
# column A (Player ID) and beyond shifts one column to the right.
$ws.Columns("A").Insert()

# Header label for the new column (row 3 is the real header row for this sheet).
$ws.Range("A3").Value = "Match ID"

# Populate the new "Match ID" column for every data row (rows 4-19) and the
# hidden totals row (row 20) with the match id value 25.
$ws.Range("A4:A20").Value = 25

# Match the bold "header" style (font 1, no border) used elsewhere on rows 3-19
# of column A, but leave row 20 (the hidden totals row) with the default style.
$ws.Range("A3:A19").Font.Bold = $true

# Restore the selection to the newly added column, as the author did after
# adding it.
[void]$ws.Range("A3:A19").Select()
